# Insert a new weekly record at row 183 for Hortaliza / Femacal de La Calera - Cilantro.
# This shifts the existing rows 183-299 down to 184-300 (dimension grows from
# A1:R299 to A1:R300) and populates the newly-opened row 183 with the new
# observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 183..299 down by one to make room for the new record.
$ws.Rows("183:183").Insert()

# Fill in the new row 183 with the new weekly observation.
$ws.Range("A183").Value = 3
$ws.Range("B183").Value = "Femacal de La Calera"
$ws.Range("C183").Value = "Coquimbo"
$ws.Range("D183").Value = 44606
$ws.Range("E183").Value = 5
$ws.Range("F183").Value = 100112040
$ws.Range("G183").Value = "Cilantro"
$ws.Range("H183").Value = "Sin especificar"
$ws.Range("I183").Value = "Primera"
$ws.Range("J183").Value = 130
$ws.Range("K183").Value = 5000
$ws.Range("L183").Value = 5500
$ws.Range("M183").Value = 5231
$ws.Range("N183").Value = "$/docena de atados (3 kilos)"
$ws.Range("O183").Value = "Provincia de Quillota"
$ws.Range("P183").Value = 1744
$ws.Range("Q183").Value = 3
$ws.Range("R183").Value = "Hortaliza"
